# This script reproduces a weekly re-shuffle of the "Cilantro" daily price
# observations: columns Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Origen, Precio $/Kg and Kg o Unidades are
# redistributed across the existing data rows (2-43). All other columns
# (Mercado ID, Mercado, Region, Codreg, Categoria ID/Categoria, Variedad,
# Calidad, Clasificacion) are left untouched, as is each row's own style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# row 2 now takes its Fecha..Kg values from former row 27
$ws.Range("D2").Value = 44272
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1893
$ws.Range("N2").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O2").Value = 'Provincia de Diguillín'
$ws.Range("P2").Value = 1893
$ws.Range("Q2").Value = 1

# row 3 now takes its Fecha..Kg values from former row 38
$ws.Range("D3").Value = 44663
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 550
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = 575
$ws.Range("N3").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O3").Value = 'Provincia de Diguillín'
$ws.Range("P3").Value = 575
$ws.Range("Q3").Value = 1

# row 4 now takes its Fecha..Kg values from former row 19
$ws.Range("D4").Value = 44671
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 550
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 575
$ws.Range("N4").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O4").Value = 'Provincia de Diguillín'
$ws.Range("P4").Value = 575
$ws.Range("Q4").Value = 1

# row 5 now takes its Fecha..Kg values from former row 41
$ws.Range("D5").Value = 44211
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1883
$ws.Range("N5").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O5").Value = 'Provincia de Diguillín'
$ws.Range("P5").Value = 1883
$ws.Range("Q5").Value = 1

# row 6 now takes its Fecha..Kg values from former row 39
$ws.Range("D6").Value = 44524
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("N6").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O6").Value = 'Provincia de Diguillín'
$ws.Range("P6").Value = 2000
$ws.Range("Q6").Value = 1

# row 7 now takes its Fecha..Kg values from former row 8
$ws.Range("D7").Value = 44539
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2200
$ws.Range("M7").Value = 2100
$ws.Range("N7").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O7").Value = 'Provincia de Diguillín'
$ws.Range("P7").Value = 2100
$ws.Range("Q7").Value = 1

# row 8 now takes its Fecha..Kg values from former row 24
$ws.Range("D8").Value = 44630
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 550
$ws.Range("L8").Value = 600
$ws.Range("M8").Value = 575
$ws.Range("N8").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O8").Value = 'Provincia de Diguillín'
$ws.Range("P8").Value = 575
$ws.Range("Q8").Value = 1

# row 9 now takes its Fecha..Kg values from former row 14
$ws.Range("D9").Value = 44166
$ws.Range("J9").Value = 240
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = 633
$ws.Range("N9").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O9").Value = 'Provincia de Diguillín'
$ws.Range("P9").Value = 633
$ws.Range("Q9").Value = 1

# row 10 now takes its Fecha..Kg values from former row 20
$ws.Range("D10").Value = 44273
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 1800
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1914
$ws.Range("N10").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O10").Value = 'Provincia de Diguillín'
$ws.Range("P10").Value = 1914
$ws.Range("Q10").Value = 1

# row 11 now takes its Fecha..Kg values from former row 2
$ws.Range("D11").Value = 44263
$ws.Range("J11").Value = 140
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1914
$ws.Range("N11").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O11").Value = 'Provincia de Diguillín'
$ws.Range("P11").Value = 1914
$ws.Range("Q11").Value = 1

# row 12 now takes its Fecha..Kg values from former row 40
$ws.Range("D12").Value = 44620
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 550
$ws.Range("L12").Value = 600
$ws.Range("M12").Value = 575
$ws.Range("N12").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O12").Value = 'Provincia de Diguillín'
$ws.Range("P12").Value = 575
$ws.Range("Q12").Value = 1

# row 13 now takes its Fecha..Kg values from former row 5
$ws.Range("D13").Value = 44208
$ws.Range("J13").Value = 130
$ws.Range("K13").Value = 1800
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 1908
$ws.Range("N13").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O13").Value = 'Provincia de Cautín'
$ws.Range("P13").Value = 1908
$ws.Range("Q13").Value = 1

# row 14 now takes its Fecha..Kg values from former row 29
$ws.Range("D14").Value = 44608
$ws.Range("J14").Value = 160
$ws.Range("K14").Value = 550
$ws.Range("L14").Value = 600
$ws.Range("M14").Value = 575
$ws.Range("N14").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O14").Value = 'Provincia de Diguillín'
$ws.Range("P14").Value = 575
$ws.Range("Q14").Value = 1

# row 15 now takes its Fecha..Kg values from former row 18
$ws.Range("D15").Value = 44270
$ws.Range("J15").Value = 260
$ws.Range("K15").Value = 1800
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 1908
$ws.Range("N15").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O15").Value = 'Provincia de Diguillín'
$ws.Range("P15").Value = 1908
$ws.Range("Q15").Value = 1

# row 16 now takes its Fecha..Kg values from former row 35
$ws.Range("D16").Value = 44265
$ws.Range("J16").Value = 220
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 1909
$ws.Range("N16").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O16").Value = 'Provincia de Diguillín'
$ws.Range("P16").Value = 1909
$ws.Range("Q16").Value = 1

# row 17 now takes its Fecha..Kg values from former row 22
$ws.Range("D17").Value = 44623
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 550
$ws.Range("L17").Value = 600
$ws.Range("M17").Value = 575
$ws.Range("N17").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O17").Value = 'Provincia de Diguillín'
$ws.Range("P17").Value = 575
$ws.Range("Q17").Value = 1

# row 18 now takes its Fecha..Kg values from former row 15
$ws.Range("D18").Value = 44532
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 2000
$ws.Range("L18").Value = 2200
$ws.Range("M18").Value = 2100
$ws.Range("N18").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O18").Value = 'Provincia de Diguillín'
$ws.Range("P18").Value = 2100
$ws.Range("Q18").Value = 1

# row 19 now takes its Fecha..Kg values from former row 36
$ws.Range("D19").Value = 44635
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 550
$ws.Range("L19").Value = 600
$ws.Range("M19").Value = 575
$ws.Range("N19").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O19").Value = 'Provincia de Diguillín'
$ws.Range("P19").Value = 575
$ws.Range("Q19").Value = 1

# row 20 now takes its Fecha..Kg values from former row 42
$ws.Range("D20").Value = 44266
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 1800
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 1913
$ws.Range("N20").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O20").Value = 'Provincia de Diguillín'
$ws.Range("P20").Value = 1913
$ws.Range("Q20").Value = 1

# row 21 now takes its Fecha..Kg values from former row 32
$ws.Range("D21").Value = 44159
$ws.Range("J21").Value = 55
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 7455
$ws.Range("N21").Value = '$/caja 36 atados'
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 207
$ws.Range("Q21").Value = 36

# row 22 now takes its Fecha..Kg values from former row 37
$ws.Range("D22").Value = 44271
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 1800
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 1920
$ws.Range("N22").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O22").Value = 'Provincia de Diguillín'
$ws.Range("P22").Value = 1920
$ws.Range("Q22").Value = 1

# row 23 now takes its Fecha..Kg values from former row 4
$ws.Range("D23").Value = 44670
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 550
$ws.Range("L23").Value = 600
$ws.Range("M23").Value = 575
$ws.Range("N23").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O23").Value = 'Provincia de Diguillín'
$ws.Range("P23").Value = 575
$ws.Range("Q23").Value = 1

# row 24 now takes its Fecha..Kg values from former row 16
$ws.Range("D24").Value = 44665
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 550
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 575
$ws.Range("N24").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O24").Value = 'Provincia de Diguillín'
$ws.Range("P24").Value = 575
$ws.Range("Q24").Value = 1

# row 25 now takes its Fecha..Kg values from former row 17
$ws.Range("D25").Value = 44610
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 550
$ws.Range("L25").Value = 600
$ws.Range("M25").Value = 575
$ws.Range("N25").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O25").Value = 'Provincia de Diguillín'
$ws.Range("P25").Value = 575
$ws.Range("Q25").Value = 1

# row 26 now takes its Fecha..Kg values from former row 13
$ws.Range("D26").Value = 44536
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = 2000
$ws.Range("N26").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O26").Value = 'Provincia de Diguillín'
$ws.Range("P26").Value = 2000
$ws.Range("Q26").Value = 1

# row 27 now takes its Fecha..Kg values from former row 3
$ws.Range("D27").Value = 44260
$ws.Range("J27").Value = 220
$ws.Range("K27").Value = 1800
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = 1909
$ws.Range("N27").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O27").Value = 'Provincia de Diguillín'
$ws.Range("P27").Value = 1909
$ws.Range("Q27").Value = 1

# row 28 now takes its Fecha..Kg values from former row 43
$ws.Range("D28").Value = 44609
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 550
$ws.Range("L28").Value = 600
$ws.Range("M28").Value = 575
$ws.Range("N28").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O28").Value = 'Provincia de Diguillín'
$ws.Range("P28").Value = 575
$ws.Range("Q28").Value = 1

# row 29 now takes its Fecha..Kg values from former row 23
$ws.Range("D29").Value = 44656
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 650
$ws.Range("M29").Value = 625
$ws.Range("N29").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O29").Value = 'Provincia de Diguillín'
$ws.Range("P29").Value = 625
$ws.Range("Q29").Value = 1

# row 30 now takes its Fecha..Kg values from former row 7
$ws.Range("D30").Value = 44631
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 550
$ws.Range("L30").Value = 600
$ws.Range("M30").Value = 575
$ws.Range("N30").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O30").Value = 'Provincia de Diguillín'
$ws.Range("P30").Value = 575
$ws.Range("Q30").Value = 1

# row 31 now takes its Fecha..Kg values from former row 6
$ws.Range("D31").Value = 44264
$ws.Range("J31").Value = 130
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = 1908
$ws.Range("N31").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O31").Value = 'Provincia de Diguillín'
$ws.Range("P31").Value = 1908
$ws.Range("Q31").Value = 1

# row 32 now takes its Fecha..Kg values from former row 12
$ws.Range("D32").Value = 44637
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 550
$ws.Range("L32").Value = 600
$ws.Range("M32").Value = 575
$ws.Range("N32").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O32").Value = 'Provincia de Diguillín'
$ws.Range("P32").Value = 575
$ws.Range("Q32").Value = 1

# row 33 now takes its Fecha..Kg values from former row 26
$ws.Range("D33").Value = 44649
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 650
$ws.Range("M33").Value = 625
$ws.Range("N33").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O33").Value = 'Provincia de Diguillín'
$ws.Range("P33").Value = 625
$ws.Range("Q33").Value = 1

# row 34 now takes its Fecha..Kg values from former row 21
$ws.Range("D34").Value = 44664
$ws.Range("J34").Value = 200
$ws.Range("K34").Value = 550
$ws.Range("L34").Value = 600
$ws.Range("M34").Value = 575
$ws.Range("N34").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O34").Value = 'Provincia de Diguillín'
$ws.Range("P34").Value = 575
$ws.Range("Q34").Value = 1

# row 35 now takes its Fecha..Kg values from former row 25
$ws.Range("D35").Value = 44267
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 1913
$ws.Range("N35").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O35").Value = 'Provincia de Diguillín'
$ws.Range("P35").Value = 1913
$ws.Range("Q35").Value = 1

# row 36 now takes its Fecha..Kg values from former row 11
$ws.Range("D36").Value = 44525
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 2000
$ws.Range("N36").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O36").Value = 'Provincia de Diguillín'
$ws.Range("P36").Value = 2000
$ws.Range("Q36").Value = 1

# row 37 now takes its Fecha..Kg values from former row 28
$ws.Range("D37").Value = 44533
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 2200
$ws.Range("M37").Value = 2100
$ws.Range("N37").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O37").Value = 'Provincia de Diguillín'
$ws.Range("P37").Value = 2100
$ws.Range("Q37").Value = 1

# row 38 now takes its Fecha..Kg values from former row 34
$ws.Range("D38").Value = 44628
$ws.Range("J38").Value = 240
$ws.Range("K38").Value = 550
$ws.Range("L38").Value = 600
$ws.Range("M38").Value = 575
$ws.Range("N38").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O38").Value = 'Provincia de Diguillín'
$ws.Range("P38").Value = 575
$ws.Range("Q38").Value = 1

# row 39 now takes its Fecha..Kg values from former row 10
$ws.Range("D39").Value = 44614
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 550
$ws.Range("L39").Value = 600
$ws.Range("M39").Value = 575
$ws.Range("N39").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O39").Value = 'Provincia de Diguillín'
$ws.Range("P39").Value = 575
$ws.Range("Q39").Value = 1

# row 40 now takes its Fecha..Kg values from former row 30
$ws.Range("D40").Value = 44644
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 550
$ws.Range("L40").Value = 600
$ws.Range("M40").Value = 575
$ws.Range("N40").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O40").Value = 'Provincia de Diguillín'
$ws.Range("P40").Value = 575
$ws.Range("Q40").Value = 1

# row 41 now takes its Fecha..Kg values from former row 9
$ws.Range("D41").Value = 44659
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 550
$ws.Range("L41").Value = 600
$ws.Range("M41").Value = 575
$ws.Range("N41").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O41").Value = 'Provincia de Diguillín'
$ws.Range("P41").Value = 575
$ws.Range("Q41").Value = 1

# row 42 now takes its Fecha..Kg values from former row 33
$ws.Range("D42").Value = 44651
$ws.Range("J42").Value = 200
$ws.Range("K42").Value = 600
$ws.Range("L42").Value = 650
$ws.Range("M42").Value = 625
$ws.Range("N42").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O42").Value = 'Provincia de Diguillín'
$ws.Range("P42").Value = 625
$ws.Range("Q42").Value = 1

# row 43 now takes its Fecha..Kg values from former row 31
$ws.Range("D43").Value = 44160
$ws.Range("J43").Value = 190
$ws.Range("K43").Value = 1300
$ws.Range("L43").Value = 1500
$ws.Range("M43").Value = 1395
$ws.Range("N43").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O43").Value = 'Provincia de Diguillín'
$ws.Range("P43").Value = 930
$ws.Range("Q43").Value = 1.5
